# Apply cryptocurrency market data updates to Sheet1
# Generated from the authoritative diff of the workbook XML
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "57.159.01"
$ws.Range("E2").Value = "  +4.61%  "
$ws.Range("D3").Value = "2.332.96"
$ws.Range("E3").Value = "  +2.09%  "
$ws.Range("E4").Value = "  +0.36%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "519.78"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +3.19%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "135.24"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +3.95%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.999"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +0.27%  "
$ws.Range("E8").Value = "  +2.06%  "
$ws.Range("D9").Value = "2.363.92"
$ws.Range("E9").Value = "  +3.00%  "
$ws.Range("E10").Value = "  +8.09%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.154"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +0.90%  "
$ws.Range("E12").Value = "  +6.63%  "
$ws.Range("E13").Value = "  +1.69%  "
$ws.Range("B14").Value = "WrappedliquidstakedEther2.0"
$ws.Range("C14").Value = "https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth"
$ws.Range("D14").Value = "2.788.49"
$ws.Range("E14").Value = "  +3.60%  "
$ws.Range("B15").Value = "Avalanche"
$ws.Range("C15").Value = "https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "23.91"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +3.24%  "
$ws.Range("D16").Value = "57.055.88"
$ws.Range("E16").Value = "  +4.39%  "
$ws.Range("E17").Value = "  +3.61%  "
$ws.Range("D18").Value = "2.364.42"
$ws.Range("E18").Value = "  +2.69%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "10.56"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +2.39%  "
$ws.Range("E20").Value = "  +3.01%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "323.12"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +5.38%  "
$ws.Range("E22").Value = "  +6.01%  "
$ws.Range("E23").Value = "  +0.61%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "61.49"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +1.17%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.998"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +0.34%  "
$ws.Range("E26").Value = "  +6.69%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "7.80"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +5.32%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "170.69"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -0.31%  "
$ws.Range("D29").Value = "0.0₃0742"
$ws.Range("E29").Value = "  +4.72%  "
$ws.Range("E30").Value = "  +9.36%  "
$ws.Range("E31").Value = "  +4.17%  "
$ws.Range("E32").Value = "  +3.58%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "18.40"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +2.28%  "
$ws.Range("E34").Value = "  +0.08%  "
$ws.Range("B35").Value = "SuiNetwork"
$ws.Range("C35").Value = "https://coinranking.com/coin/3xJluUMvp+suinetwork-sui"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.956"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +1.83%  "
$ws.Range("B36").Value = "FirstDigitalUSD"
$ws.Range("C36").Value = "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.997"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +0.23%  "
$ws.Range("E37").Value = "  +5.08%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "4.02"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +7.13%  "
$ws.Range("E39").Value = "  +7.21%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "37.55"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +3.70%  "
$ws.Range("E41").Value = "  +1.88%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "140.67"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +12.12%  "
$ws.Range("E43").Value = "  +5.51%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "278.55"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +12.87%  "
$ws.Range("E45").Value = "  +1.77%  "
$ws.Range("E46").Value = "  +3.61%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.0931"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +3.52%  "
$ws.Range("E49").Value = "  +2.09%  "
$ws.Range("E50").Value = "  +4.39%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "17.02"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +3.22%  "
